$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1764705882352941
$ws.Range("C2").Value = 0.5640138408304498
$ws.Range("J2").Value = 0.01384083044982699
$ws.Range("P2").Value = 0.1591695501730104
$ws.Range("S2").Value = 0.08650519031141868
$ws.Range("B3").Value = 0.01129943502824859
$ws.Range("C3").Value = 0.01129943502824859
$ws.Range("J3").Value = 0.03954802259887006
$ws.Range("P3").Value = 0.768361581920904
$ws.Range("S3").Value = 0.1694915254237288
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.2333333333333333
$ws.Range("B6").Value = 0.06486486486486487
$ws.Range("D6").Value = 0.01621621621621622
$ws.Range("F6").Value = 0.06486486486486487
$ws.Range("J6").Value = 0.2648648648648649
$ws.Range("O6").Value = 0.01081081081081081
$ws.Range("Q6").Value = 0.145945945945946
$ws.Range("R6").Value = 0.04864864864864865
$ws.Range("S6").Value = 0.3837837837837838
$ws.Range("B7").Value = 0.1414141414141414
$ws.Range("D7").Value = 0.0101010101010101
$ws.Range("F7").Value = 0.06060606060606061
$ws.Range("J7").Value = 0.1363636363636364
$ws.Range("O7").Value = 0.03535353535353535
$ws.Range("Q7").Value = 0.1363636363636364
$ws.Range("R7").Value = 0.08585858585858586
$ws.Range("S7").Value = 0.3939393939393939
$ws.Range("B8").Value = 0.09424083769633508
$ws.Range("D8").Value = 0.01047120418848168
$ws.Range("F8").Value = 0.07591623036649214
$ws.Range("J8").Value = 0.1282722513089005
$ws.Range("O8").Value = 0.02879581151832461
$ws.Range("Q8").Value = 0.1544502617801047
$ws.Range("R8").Value = 0.07329842931937172
$ws.Range("S8").Value = 0.4345549738219895
$ws.Range("B9").Value = 0.08641975308641975
$ws.Range("D9").Value = 0.01234567901234568
$ws.Range("E9").Value = 0.006172839506172839
$ws.Range("F9").Value = 0.06172839506172839
$ws.Range("J9").Value = 0.1234567901234568
$ws.Range("O9").Value = 0.02469135802469136
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.06172839506172839
$ws.Range("S9").Value = 0.4567901234567901
$ws.Range("B10").Value = 0.120704845814978
$ws.Range("D10").Value = 0.01585903083700441
$ws.Range("F10").Value = 0.05286343612334802
$ws.Range("J10").Value = 0.1427312775330397
$ws.Range("O10").Value = 0.01938325991189427
$ws.Range("Q10").Value = 0.2237885462555066
$ws.Range("R10").Value = 0.06167400881057269
$ws.Range("S10").Value = 0.3629955947136564
$ws.Range("G11").Value = 0.1904761904761905
$ws.Range("J11").Value = 0.06349206349206349
$ws.Range("K11").Value = 0.2476190476190476
$ws.Range("L11").Value = 0.4761904761904762
$ws.Range("S11").Value = 0.02222222222222222
$ws.Range("G12").Value = 0.7058823529411765
$ws.Range("J12").Value = 0.261437908496732
$ws.Range("K12").Value = 0.0130718954248366
$ws.Range("L12").Value = 0.006535947712418301
$ws.Range("S12").Value = 0.0130718954248366
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.2708333333333333
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.00558659217877095
$ws.Range("H15").Value = 0.106145251396648
$ws.Range("I15").Value = 0.05586592178770949
$ws.Range("J15").Value = 0.4189944134078212
$ws.Range("K15").Value = 0.0446927374301676
$ws.Range("M15").Value = 0.0223463687150838
$ws.Range("O15").Value = 0.09497206703910614
$ws.Range("S15").Value = 0.2513966480446927
$ws.Range("F16").Value = 0.01587301587301587
$ws.Range("H16").Value = 0.1746031746031746
$ws.Range("I16").Value = 0.0582010582010582
$ws.Range("J16").Value = 0.3756613756613756
$ws.Range("K16").Value = 0.1322751322751323
$ws.Range("M16").Value = 0.03174603174603174
$ws.Range("O16").Value = 0.06878306878306878
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.02267002518891688
$ws.Range("H17").Value = 0.1486146095717884
$ws.Range("I17").Value = 0.1007556675062972
$ws.Range("J17").Value = 0.4483627204030227
$ws.Range("K17").Value = 0.09823677581863979
$ws.Range("M17").Value = 0.01007556675062972
$ws.Range("O17").Value = 0.03778337531486146
$ws.Range("S17").Value = 0.1335012594458438
$ws.Range("F18").Value = 0.02205882352941177
$ws.Range("H18").Value = 0.1764705882352941
$ws.Range("I18").Value = 0.09558823529411764
$ws.Range("J18").Value = 0.4558823529411765
$ws.Range("K18").Value = 0.08823529411764706
$ws.Range("M18").Value = 0.02941176470588235
$ws.Range("O18").Value = 0.04411764705882353
$ws.Range("S18").Value = 0.08823529411764706
$ws.Range("F19").Value = 0.02042628774422735
$ws.Range("H19").Value = 0.2238010657193606
$ws.Range("J19").Value = 0.3374777975133215
$ws.Range("K19").Value = 0.1385435168738899
$ws.Range("M19").Value = 0.02753108348134991
$ws.Range("N19").Value = 0.001776198934280639
$ws.Range("O19").Value = 0.0541740674955595
$ws.Range("S19").Value = 0.1181172291296625
$ws.Range("I19").Value = 0.07815275310834814
